# New crime data collected
# Updates the weekly CompStat report: header volume/date strings, and the
# weekly/28-day/YTD/2-year crime figures in rows 16-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: "Volume 30   Number  45" -> "...46" and the week-covering date
# range "11/6/2023 ... 11/12/2023" -> "11/13/2023 ... 11/19/2023".
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# ---------------------------------------------------------------------------
# Crime-figures table (rows 16-30). Most cells simply get a new numeric
# value while keeping their existing number-format/style.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 186
$ws.Range("J16").Value = 219
$ws.Range("K16").Value = -15.068493150684
$ws.Range("L16").Value = 18.471337579617
$ws.Range("M16").Value = 55
$ws.Range("N16").Value = -84.564315352697

$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 83.333333333333
$ws.Range("I17").Value = 170
$ws.Range("J17").Value = 152
$ws.Range("K17").Value = 11.842105263157
$ws.Range("L17").Value = 21.428571428571
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -35.849056603773

$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -4.545454545454
$ws.Range("I18").Value = 236
$ws.Range("J18").Value = 236
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 15.121951219512
$ws.Range("M18").Value = 8.256880733944
$ws.Range("N18").Value = -91.057218643425

$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = 43.478260869565
$ws.Range("F19").Value = 127
$ws.Range("G19").Value = 133
$ws.Range("H19").Value = -4.511278195488
$ws.Range("I19").Value = 1506
$ws.Range("J19").Value = 1560
$ws.Range("K19").Value = -3.461538461538
$ws.Range("L19").Value = 44.668587896253
$ws.Range("M19").Value = 34.825425246195
$ws.Range("N19").Value = -54.706766917293

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -46.666666666666
$ws.Range("I20").Value = 141
$ws.Range("J20").Value = 161
$ws.Range("K20").Value = -12.422360248447
$ws.Range("L20").Value = -7.843137254901
$ws.Range("M20").Value = 69.879518072289
$ws.Range("N20").Value = -95.411649853563

$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 184
$ws.Range("G21").Value = 192
$ws.Range("H21").Value = -4.166666666666
$ws.Range("I21").Value = 2254
$ws.Range("J21").Value = 2340
$ws.Range("K21").Value = -3.675213675213
$ws.Range("L21").Value = 31.812865497076
$ws.Range("M21").Value = 37.690897984117
$ws.Range("N21").Value = -78.616829522815

# Row 22 (Transit): C22 flips from the text placeholder "0" to the number 1,
# so also restore the numeric (#,##0) style used by its neighbours.
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("I22").Value = 41
$ws.Range("J22").Value = 42
$ws.Range("K22").Value = -2.380952380952
$ws.Range("L22").Value = 192.857142857143
$ws.Range("M22").Value = 70.833333333333

# Row 23 (Housing): C23 flips from the number 3 to the text placeholder "0".
# Build that text value in a scratch cell (quote-prefixed so it is stored as
# text, not auto-converted to a number), copy/paste its format onto C23, and
# finally stamp C23 with the same right-aligned "General" style used by the
# other text-placeholder cells (taken from A13-style cell C14).
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("J23").Value = 33
$ws.Range("K23").Value = -21.212121212121

$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 101
$ws.Range("E24").Value = -60.396039603960
$ws.Range("F24").Value = 200
$ws.Range("G24").Value = 366
$ws.Range("H24").Value = -45.355191256830
$ws.Range("I24").Value = 2768
$ws.Range("J24").Value = 3588
$ws.Range("K24").Value = -22.853957636566
$ws.Range("L24").Value = 27.264367816092
$ws.Range("M24").Value = 83.798140770252

$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 51.851851851851
$ws.Range("I25").Value = 317
$ws.Range("J25").Value = 334
$ws.Range("K25").Value = -5.089820359281
$ws.Range("L25").Value = 3.934426229508
$ws.Range("M25").Value = 1.277955271565

$ws.Range("L26").Value = 15

# Row 27 (Other Sex Crimes): C27 flips from the number 1 to the text
# placeholder "0", same treatment as C23 above.
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = -42.857142857142
$ws.Range("J27").Value = 73
$ws.Range("K27").Value = 5.479452054794
$ws.Range("L27").Value = -9.411764705882

# Row 28 (Shooting Vic.): F28 flips from the number 1 to the text
# placeholder "0".
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("F28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)

# Row 29 (Shooting Inc.): F29 flips from the number 1 to the text
# placeholder "0".
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("F29").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 19
$ws.Range("K30").Value = -29.629629629629
$ws.Range("L30").Value = 35.714285714285

# Clean up the scratch cell used for the text-placeholder conversions above.
$ws.Range("ZZ1").Clear()
